$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert the header row to the pre-"Document panel loading" naming scheme.
$ws.Range("A1").Value = "HGNC_IDnumber"
$ws.Range("B1").Value = "HGNC_symbol"
$ws.Range("C1").Value = "Disease_associated_transcript"
$ws.Range("D1").Value = "Genetic_disease_model"
$ws.Range("E1").Value = "Mosaicism"
$ws.Range("F1").Value = "Reduced_penetrance"
$ws.Range("G1").Value = "Database_entry_version"

# Restore the original selection.
$ws.Range("E1").Select()
